# Apply the "/presentation/primarylanguage" documentation update (related to #6)
# to the first table of the document.
#
# Row 12 (1-indexed) = "/presentation/primarylanguage" row
# Row 13 (1-indexed) = "/presentation/language" row

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wdColorYellow = 65535      # RGB(255,255,0)  -> FFFF00
$wdColorGreen  = 1758337    # RGB(0x81,0xD4,0x1A) -> 81D41A
$wdColorRed    = 255        # RGB(255,0,0)    -> FF0000

function Set-CellShading($cell, $color) {
    $cell.Shading.BackgroundPatternColor = $color
    $cell.Shading.Texture = 0
}

# ---------------------------------------------------------------------
# Row 12: /presentation/primarylanguage
# ---------------------------------------------------------------------

# Column 6: "presentation_language" -> "presentation_primarylanguage"
$cell = $t.Cell(12, 6)
$cell.Range.Text = "presentation_primarylanguage"

# Column 7: empty -> two paragraphs, shading -> yellow
$cell = $t.Cell(12, 7)
$cell.Range.Text = "presentation_primarylanguage" + [char]13 + "can be wrong with single language display " + [char]8211 + " warn is displayed"
Set-CellShading $cell $wdColorYellow

# Column 8: "No response" -> "presentation_primarylanguage", shading -> green
$cell = $t.Cell(12, 8)
$cell.Range.Text = "presentation_primarylanguage"
Set-CellShading $cell $wdColorGreen

# Column 10: "#6 no functionality" -> three paragraphs, shading FF0000 -> yellow
$cell = $t.Cell(12, 10)
$cell.Range.Text = "#6 with 6.10" + [char]13 + [char]13 + "incorrect when single language display only"
Set-CellShading $cell $wdColorYellow

# ---------------------------------------------------------------------
# Row 13: /presentation/language
# ---------------------------------------------------------------------

# Column 1: "/presentation/language" -> "/presentation/languages"
$cell = $t.Cell(13, 1)
$r = $cell.Range
$appendRange = $d.Range($r.Start, $r.End - 1)
$appendRange.InsertAfter("s")

# Column 6: "presentation_language" -> "presentation_languages" + line breaks + note, shading -> red
$cell = $t.Cell(13, 6)
$cell.Range.Text = "presentation_languages" + [char]11 + [char]11 + "not setting languages"
Set-CellShading $cell $wdColorRed

# Column 7: empty -> "presentation_languages" + line breaks + note, shading -> red
$cell = $t.Cell(13, 7)
$cell.Range.Text = "presentation_languages" + [char]11 + [char]11 + "disabled because incorrect"
Set-CellShading $cell $wdColorRed

# Column 8: "No response" -> "presentation_languages" + line breaks + note, shading -> red
$cell = $t.Cell(13, 8)
$cell.Range.Text = "presentation_languages" + [char]11 + [char]11 + "disabled because incorrect"
Set-CellShading $cell $wdColorRed

# Column 10: "#6 no functionality" -> "#6 open bugs – therefore unuseable with 6.10"
$cell = $t.Cell(13, 10)
$cell.Range.Text = "#6 open bugs " + [char]8211 + " therefore unuseable with 6.10"

Write-Host "Done applying presentation/primarylanguage edits."
